# Atualização automática da planilha
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "Projetos" sheet: rename project P11 (row 13) from
#    "Contas a Receber" to "Contratação", and clear its "Área"
#    column (C13), which used to repeat the same text.
# ---------------------------------------------------------------
$wsProjetos = $wb.Worksheets.Item("Projetos")
$wsProjetos.Range("B13").Value = "Contratação"
$wsProjetos.Range("C13").Value = ""

# ---------------------------------------------------------------
# 2) "Etapas" sheet: the same project's phase rows (54:57) used to
#    read "ITSM" in column B - update them to the new project name.
#    Then append the five standard phase rows (Levantamento de
#    requisitos / Configuração dos módulos / Homologação e testes /
#    Treinamento de usuários / Go-Live e estabilização) for four
#    more project/area combinations.
# ---------------------------------------------------------------
$wsEtapas = $wb.Worksheets.Item("Etapas")

$wsEtapas.Range("B53:B57").Value = "Contratação"

# P12 / Financeiro -> rows 58:62 (same plain style as the template
# rows copied from, no extra treatment needed on column B).
$wsEtapas.Range("A48:H52").Copy($wsEtapas.Range("A58:H62"))
$wsEtapas.Range("A58:A62").Value = "P12"
$wsEtapas.Range("B58:B62").Value = "Financeiro"

# P10 / Engenharia - Entrega -> rows 63:67 (column B keeps the bold
# style used on the sheet's first block, so re-apply that format).
$wsEtapas.Range("A48:H52").Copy($wsEtapas.Range("A63:H67"))
$wsEtapas.Range("A63:A67").Value = "P10"
$wsEtapas.Range("B63:B67").Value = "Engenharia - Entrega"
$wsEtapas.Range("B3:B7").Copy()
$wsEtapas.Range("B63:B67").PasteSpecial(-4122)

# P11 / Engenharia - Execução -> rows 68:72
$wsEtapas.Range("A48:H52").Copy($wsEtapas.Range("A68:H72"))
$wsEtapas.Range("A68:A72").Value = "P11"
$wsEtapas.Range("B68:B72").Value = "Engenharia - Execução"
$wsEtapas.Range("B3:B7").Copy()
$wsEtapas.Range("B68:B72").PasteSpecial(-4122)

# P10 / Engenharia - Legalização -> rows 73:77
$wsEtapas.Range("A48:H52").Copy($wsEtapas.Range("A73:H77"))
$wsEtapas.Range("A73:A77").Value = "P10"
$wsEtapas.Range("B73:B77").Value = "Engenharia - Legalização"
$wsEtapas.Range("B3:B7").Copy()
$wsEtapas.Range("B73:B77").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 3) View state: the workbook was left with "Etapas" as the active
#    sheet/tab, scrolled down to the newly added rows.
# ---------------------------------------------------------------
$wsEtapas.Activate()
$wsEtapas.Range("B73:B77").Select()
